$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Home/Away/ML_h/ML_a values for rows 2-17 (Week 16, 2021 odds update)
$data = @(
    @("Titans",     "49ers",      "+145", "-165"),
    @("Packers",    "Browns",     "-380", "+290"),
    @("Cardinals",  "Colts",      "-120", "+100"),
    @("Falcons",    "Lions",      "-250", "+200"),
    @("Bengals",    "Ravens",     "-140", "+120"),
    @("Vikings",    "Rams",       "+135", "-155"),
    @("Patriots",   "Bills",      "-140", "+120"),
    @("Jets",       "Jaguars",    "-110", "-110"),
    @("Eagles",     "Giants",     "-475", "+350"),
    @("Panthers",   "Buccaneers", "+340", "-450"),
    @("Texans",     "Chargers",   "+350", "-475"),
    @("Seahawks",   "Bears",      "-275", "+220"),
    @("Chiefs",     "Steelers",   "-350", "+270"),
    @("Raiders",    "Broncos",    "-110", "-110"),
    @("Cowboys",    "Washington", "-475", "+350"),
    @("Saints",     "Dolphins",   "-160", "+140")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 2).Value = $entry[0]
    $ws.Cells.Item($row, 3).Value = $entry[1]
    # Prefix with an apostrophe so Excel keeps these odds (e.g. "+145", "-380")
    # as literal text instead of converting them to numbers.
    $ws.Cells.Item($row, 4).Value = "'" + $entry[2]
    $ws.Cells.Item($row, 5).Value = "'" + $entry[3]
    $row++
}
